$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tidy")

# Column D holds the "estimate" values. Rows 2-50 contain data whose sign
# needs to be flipped (corrected estimate sign).
for ($row = 2; $row -le 50; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $current = $cell.Value()
    $cell.Value = -1 * $current
}
